$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (bold, border, centered) from A12 into A13, matching the year-label column format
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 528
$ws.Range("C13").Value = 3048
$ws.Range("D13").Value = 215
$ws.Range("E13").Value = 1137
$ws.Range("F13").Value = 60
$ws.Range("G13").Value = 15
$ws.Range("H13").Value = 190
$ws.Range("I13").Value = 2592
$ws.Range("J13").Value = 32
$ws.Range("K13").Value = 568
$ws.Range("L13").Value = 13
$ws.Range("M13").Value = 4
$ws.Range("N13").Value = 947
$ws.Range("P13").Value = 4379
$ws.Range("R13").Value = 14
$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 1
$ws.Range("V13").Value = 1
$ws.Range("W13").Value = 1
$ws.Range("Y13").Value = 13
$ws.Range("Z13").Value = 2
$ws.Range("AC13").Value = 1
$ws.Range("AD13").Value = 44
$ws.Range("AE13").Value = 1
$ws.Range("AH13").Value = 5
$ws.Range("AI13").Value = 20
$ws.Range("AJ13").Value = 295
$ws.Range("AK13").Value = 423
$ws.Range("AL13").Value = 158
$ws.Range("AM13").Value = 382
$ws.Range("AN13").Value = 27
$ws.Range("AO13").Value = 9
$ws.Range("AP13").Value = 287
$ws.Range("AR13").Value = 1692
$ws.Range("AS13").Value = 21
$ws.Range("AT13").Value = 3
$ws.Range("AU13").Value = 2
$ws.Range("AV13").Value = 32
$ws.Range("AW13").Value = 38
$ws.Range("AX13").Value = 1238
$ws.Range("AY13").Value = 8
$ws.Range("AZ13").Value = 6693
$ws.Range("BA13").Value = 42
$ws.Range("BB13").Value = 32
$ws.Range("BC13").Value = 25
$ws.Range("BD13").Value = 174
$ws.Range("BE13").Value = 18
$ws.Range("BF13").Value = 2
$ws.Range("BG13").Value = 4
$ws.Range("BH13").Value = 7
$ws.Range("BI13").Value = 578
$ws.Range("BJ13").Value = 25
$ws.Range("BL13").Value = 14
$ws.Range("BM13").Value = 24
$ws.Range("BN13").Value = 210
$ws.Range("BO13").Value = 47
$ws.Range("BP13").Value = 17
$ws.Range("BQ13").Value = 16
$ws.Range("BR13").Value = 62
$ws.Range("BS13").Value = 269

# Empty-but-present text cells (mirrors the "<c t=\"inlineStr\"/>" cells used for
# not-applicable columns in this sheet): enter a lone quote-prefix so the cell
# commits as an empty Text value instead of being cleared, then drop the implicit
# quote-prefix formatting so the cell keeps the sheet default style.
$ws.Range("O13").Value = "'"
$ws.Range("O13").Style = "Normal"
$ws.Range("Q13").Value = "'"
$ws.Range("Q13").Style = "Normal"
$ws.Range("S13").Value = "'"
$ws.Range("S13").Style = "Normal"
$ws.Range("X13").Value = "'"
$ws.Range("X13").Style = "Normal"
$ws.Range("AA13").Value = "'"
$ws.Range("AA13").Style = "Normal"
$ws.Range("AB13").Value = "'"
$ws.Range("AB13").Style = "Normal"
$ws.Range("AF13").Value = "'"
$ws.Range("AF13").Style = "Normal"
$ws.Range("AG13").Value = "'"
$ws.Range("AG13").Style = "Normal"
$ws.Range("AQ13").Value = "'"
$ws.Range("AQ13").Style = "Normal"
$ws.Range("BK13").Value = "'"
$ws.Range("BK13").Style = "Normal"
